$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'" + '67.041.26'
$ws.Range('E2').Value = "'" + '  +0.13%  '
$ws.Range('D3').Value = "'" + '3.131.46'
$ws.Range('E3').Value = "'" + '  +0.94%  '
$ws.Range('E4').Value = "'" + '  -0.01%  '
$ws.Range('D5').Value = "'" + '581.31'
$ws.Range('E5').Value = "'" + '  +0.23%  '
$ws.Range('D6').Value = "'" + '173.83'
$ws.Range('E6').Value = "'" + '  +0.01%  '
$ws.Range('D7').Value = "'" + '1.00'
$ws.Range('E7').Value = "'" + '  -0.01%  '
$ws.Range('D8').Value = "'" + '0.523'
$ws.Range('E8').Value = "'" + '  -0.19%  '
$ws.Range('D9').Value = "'" + '6.44'
$ws.Range('E9').Value = "'" + '  -1.00%  '
$ws.Range('D10').Value = "'" + '0.155'
$ws.Range('E10').Value = "'" + '  -0.88%  '
$ws.Range('D11').Value = "'" + '0.481'
$ws.Range('E11').Value = "'" + '  -0.27%  '
$ws.Range('D12').Value = "'" + '0.0000249'
$ws.Range('E12').Value = "'" + '  -0.47%  '
$ws.Range('D13').Value = "'" + '37.66'
$ws.Range('E13').Value = "'" + '  +1.09%  '
$ws.Range('D14').Value = "'" + '0.123'
$ws.Range('E14').Value = "'" + '  -1.50%  '
$ws.Range('D15').Value = "'" + '67.013.18'
$ws.Range('E15').Value = "'" + '  +0.07%  '
$ws.Range('D16').Value = "'" + '7.15'
$ws.Range('E16').Value = "'" + '  -0.63%  '
$ws.Range('D17').Value = "'" + '3.129.82'
$ws.Range('E17').Value = "'" + '  +0.92%  '
$ws.Range('D18').Value = "'" + '16.43'
$ws.Range('E18').Value = "'" + '  +1.41%  '
$ws.Range('D19').Value = "'" + '490.58'
$ws.Range('E19').Value = "'" + '  +1.88%  '
$ws.Range('D20').Value = "'" + '0.710'
$ws.Range('E20').Value = "'" + '  -0.90%  '
$ws.Range('D21').Value = "'" + '7.92'
$ws.Range('E21').Value = "'" + '  +5.35%  '
$ws.Range('D22').Value = "'" + '84.29'
$ws.Range('E22').Value = "'" + '  +0.15%  '
$ws.Range('D23').Value = "'" + '13.28'
$ws.Range('E23').Value = "'" + '  +1.83%  '
$ws.Range('D24').Value = "'" + '2.30'
$ws.Range('E24').Value = "'" + '  -2.17%  '
$ws.Range('D25').Value = "'" + '10.34'
$ws.Range('E25').Value = "'" + '  +3.37%  '
$ws.Range('E26').Value = "'" + '  +0.09%  '
$ws.Range('D27').Value = "'" + '7.96'
$ws.Range('E27').Value = "'" + '  -0.65%  '
$ws.Range('D28').Value = "'" + '2.37'
$ws.Range('E28').Value = "'" + '  -0.78%  '
$ws.Range('D29').Value = "'" + '2.69'
$ws.Range('E29').Value = "'" + '  +0.40%  '
$ws.Range('D30').Value = "'" + '28.74'
$ws.Range('E30').Value = "'" + '  -0.01%  '
$ws.Range('E31').Value = "'" + '  -0.49%  '
$ws.Range('E32').Value = "'" + '  -6.61%  '
$ws.Range('E33').Value = "'" + '  -0.05%  '
$ws.Range('D34').Value = "'" + '5.89'
$ws.Range('E34').Value = "'" + '  -0.13%  '
$ws.Range('D35').Value = "'" + '0.979'
$ws.Range('E35').Value = "'" + '  -2.50%  '
$ws.Range('D36').Value = "'" + '46.99'
$ws.Range('E36').Value = "'" + '  -1.50%  '
$ws.Range('D37').Value = "'" + '50.19'
$ws.Range('E37').Value = "'" + '  +0.01%  '
$ws.Range('D38').Value = "'" + '2.06'
$ws.Range('E38').Value = "'" + '  -3.40%  '
$ws.Range('D39').Value = "'" + '0.313'
$ws.Range('E39').Value = "'" + '  -1.21%  '
$ws.Range('E40').Value = "'" + '  +1.91%  '
$ws.Range('D41').Value = "'" + '8.57'
$ws.Range('E41').Value = "'" + '  -1.25%  '
$ws.Range('B42').Value = "'" + 'Bittensor'
$ws.Range('C42').Value = "'" + 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D42').Value = "'" + '386.79'
$ws.Range('E42').Value = "'" + '  +1.59%  '
$ws.Range('B43').Value = "'" + 'Maker'
$ws.Range('C43').Value = "'" + 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D43').Value = "'" + '2.825.33'
$ws.Range('E43').Value = "'" + '  +0.19%  '
$ws.Range('D44').Value = "'" + '2.60'
$ws.Range('E44').Value = "'" + '  -7.54%  '
$ws.Range('D45').Value = "'" + '0.0354'
$ws.Range('E45').Value = "'" + '  -1.85%  '
$ws.Range('D46').Value = "'" + '135.96'
$ws.Range('E46').Value = "'" + '  +0.73%  '
$ws.Range('E47').Value = "'" + '  +0.02%  '
$ws.Range('D48').Value = "'" + '25.18'
$ws.Range('E48').Value = "'" + '  +1.07%  '
$ws.Range('D49').Value = "'" + '2.23'
$ws.Range('E49').Value = "'" + '  +0.21%  '
$ws.Range('E50').Value = "'" + '  -0.19%  '
$ws.Range('D51').Value = "'" + '6.79'
$ws.Range('E51').Value = "'" + '  -0.51%  '
